$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Area / Atotal columns, plus a small J:K summary block
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# D column (segment midpoint) becomes a shared formula across D3:D9
$ws.Range("D3:D9").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# New Area column: per-segment incremental area, and a cumulative total
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Summary block mirroring the totals
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Match the saved selection state
$ws.Range("J2:K2").Select()
